# Auto-generated from the cryptos.xlsx price-refresh diff (GitHub Actions run).
# Updates Price (D) / Volume(1h) (E) text cells for rows 2-51, plus the
# Bittensor/Kaspa row swap at rows 28-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.451.29"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.629.51"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "593.46"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "152.23"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("D10").Value = "0.396"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("D11").Value = "'5.80"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "28.59"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "3.102.53"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +13.53%  "
$ws.Range("D16").Value = "64.413.45"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "2.643.02"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "12.24"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "350.09"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "7.11"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "67.46"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "1.71"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "9.28"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "1.65"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "8.27"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "543.15"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "0.973"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").Value = "0.0₃0914"
$ws.Range("E31").Value = "  +8.72%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("D34").Value = "5.58"
$ws.Range("E34").Value = "  +6.67%  "
$ws.Range("D35").Value = "6.19"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").Value = "164.03"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "20.12"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "168.29"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "41.66"
$ws.Range("E43").Value = "  +4.70%  "
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").Value = "23.43"
$ws.Range("E45").Value = "  +9.03%  "
$ws.Range("D46").Value = "0.0593"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "2.22"
$ws.Range("E47").Value = "  +12.44%  "
$ws.Range("D48").Value = "0.643"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").Value = "0.0252"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "0.0979"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'19.40"
$ws.Range("E51").Value = "  +0.57%  "
